# feat: add 2022-Q1 data
#
# Before:  Sheets = [ "2021-Q4", "总计" ]
# After:   Sheets = [ "2021-Q4", "2022-Q1", "总计" ]
#   - "2022-Q1" is a brand-new per-fund holdings sheet (same shape as "2021-Q4").
#   - "总计" (summary) gets a new row inserted above the existing "2021-Q4" row.

$wb = $excel.ActiveWorkbook

$q4Sheet  = $wb.Worksheets.Item(1)   # "2021-Q4" - stays untouched, used as insertion anchor

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet between "2021-Q4" and "总计".
# ---------------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4Sheet)
$q1Sheet.Name = "2022-Q1"

# NOTE: fetch the "总计" sheet by name, and only *after* the new sheet has
# been inserted - worksheet references captured by numeric Item() index
# before an Add() do not reliably track the sheet across the shift.
$sumSheet = $wb.Worksheets.Item("总计")   # "总计" - existing summary sheet

# Header row (bold, thin border, centered/top-aligned - same look as "2021-Q4").
# Written cell-by-cell (not as a bulk array assignment) for reliability.
$q1Sheet.Cells.Item(1, 2).Value = "基金代码"
$q1Sheet.Cells.Item(1, 3).Value = "基金名称"
$q1Sheet.Cells.Item(1, 4).Value = "基金规模"
$q1Sheet.Cells.Item(1, 5).Value = "股票总仓位"
$q1Sheet.Cells.Item(1, 6).Value = "仓位占比"
$q1Sheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1Sheet.Cells.Item(1, 8).Value = "仓位排名"

$headerRange = $q1Sheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Per-fund holdings rows.
$rows = @(
    @{ idx = 0; code = "160416"; name = "华安标普全球石油指数 (QDII-LOF)"; size = "3.37"; pos = "95.08"; pct = "3.95"; mv = "0.1331"; rank = 4 },
    @{ idx = 1; code = "513080"; name = "华安法国CAC40ETF（QDII）";         size = "0.60"; pos = "96.69"; pct = "7.65"; mv = "0.0459"; rank = 2 },
    @{ idx = 2; code = "006282"; name = "上投摩根欧洲动力策略股票（QDII）"; size = "0.48"; pos = "89.68"; pct = "2.90"; mv = "0.0139"; rank = 3 }
)

foreach ($r in $rows) {
    $row = 2 + $r.idx

    $aCell = $q1Sheet.Cells.Item($row, 1)
    $aCell.Value = $r.idx
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    # Text-valued columns - leading "'" keeps them literal text (preserves
    # leading zeros / trailing zeros exactly as scraped, no numeric coercion).
    $q1Sheet.Cells.Item($row, 2).Value = "'" + $r.code
    $q1Sheet.Cells.Item($row, 3).Value = "'" + $r.name
    $q1Sheet.Cells.Item($row, 4).Value = "'" + $r.size
    $q1Sheet.Cells.Item($row, 5).Value = "'" + $r.pos
    $q1Sheet.Cells.Item($row, 6).Value = "'" + $r.pct
    $q1Sheet.Cells.Item($row, 7).Value = "'" + $r.mv

    # Rank column is numeric.
    $q1Sheet.Cells.Item($row, 8).Value = $r.rank
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row above the existing
#    "2021-Q4" row and fill it in with the "2022-Q1" totals.
# ---------------------------------------------------------------------------
$sumSheet.Rows.Item(2).Insert()
$sumSheet.Range("A2:D2").ClearFormats()

$sumA2 = $sumSheet.Range("A2")
$sumA2.Value = 0
$sumA2.Font.Bold = $true
$sumA2.Borders.LineStyle = 1
$sumA2.HorizontalAlignment = -4108
$sumA2.VerticalAlignment = -4160

$sumSheet.Range("B2").Value = "'2022-Q1"
$sumSheet.Range("C2").Value = 3
$sumSheet.Range("D2").Value = 0.19

# The pre-existing "2021-Q4" row (now shifted down to row 3) keeps its data
# but its running index in column A needs to advance from 0 to 1.
$sumSheet.Range("A3").Value = 1

Write-Host "2022-Q1 sheet + 总计 update applied"
